$wb = $excel.ActiveWorkbook

# Work on the "Metadata" worksheet explicitly (it is the active sheet already,
# but we select it by name to be safe/robust).
$ws = $wb.Worksheets.Item("Metadata")

# Row 4 = "Name" property, column B currently empty -> set to "OrdreVs"
$ws.Range("B4").Value = "OrdreVs"

# Row 8 = "Date" property -> update the generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
